$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Hartmut"

$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

$ws.Range("D5").Value = "KONTOSTAND AM 06.02.2024"

$ws.Range("B6").Value = "08.02."
$ws.Range("C6").Value = "09.02."
$ws.Range("D6").Value = "BEITRAG Allianz SE K-90207280"
$ws.Range("E6").Value = "57,46-"

$ws.Range("B7").Value = "12.02."
$ws.Range("C7").Value = "13.02."
$ws.Range("D7").Value = "ZALANDO MKTPLC EU ZKVXCF"
$ws.Range("E7").Value = "124,75-"

$ws.Range("B8").Value = "13.02."
$ws.Range("C8").Value = "14.02."
$ws.Range("D8").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E8").Value = "25,46-"

$ws.Range("D12").Value = "KONTOSTAND AM 18.02.2024"
$ws.Range("E12").Value = "207,67-"

$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 26.02.2024"
